$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 42: Min Cost Climbing Stairs
$ws.Range("A42").Value = "Min Cost Climbing Stairs"
$ws.Range("B42").Value = "Dynamic Programming"
$ws.Range("C42").Value = "No"
$ws.Range("D42").Value = "Yes"
$ws.Range("E42").Value = "Easy"
$ws.Range("F42").Value = "Medium"
$ws.Range("G42").Value = "746 - Min Cost Climbing Stairs"
$ws.Hyperlinks.Add($ws.Range("G42"), "746 - Min Cost Climbing Stairs", "", "", "746 - Min Cost Climbing Stairs") | Out-Null
$ws.Range("G42").Style = "Hyperlink"

# Row 43: Length of Last Word
$ws.Range("A43").Value = "Length of Last Word"
$ws.Range("B43").Value = "String"
$ws.Range("C43").Value = "No"
$ws.Range("D43").Value = "No"
$ws.Range("E43").Value = "Easy"
$ws.Range("F43").Value = "Easy"
$ws.Range("G43").Value = "58 - Length of Last Word"
$ws.Hyperlinks.Add($ws.Range("G43"), "58 - Length of Last Word", "", "", "58 - Length of Last Word") | Out-Null
$ws.Range("G43").Style = "Hyperlink"

# Extend conditional formatting and data validation ranges to include new rows
$ws.Range("D2:G8,D9:F43").FormatConditions.Delete()
$cf1 = $ws.Range("D2:G8,D9:F43").FormatConditions.Add(2, 3, '"Hard"')
$cf1.Interior.ThemeColor = 9
$cf1.Priority = 7
$cf2 = $ws.Range("D2:G8,D9:F43").FormatConditions.Add(2, 3, '"Medium"')
$cf2.Interior.ThemeColor = 5
$cf2.Priority = 8
$cf3 = $ws.Range("D2:G8,D9:F43").FormatConditions.Add(2, 3, '"Easy"')
$cf3.Interior.Color = 255
$cf3.Priority = 9

$ws.Range("E2:F43").Validation.Delete()
$ws.Range("E2:F43").Validation.Add(3, 1, 1, "Easy, Medium, Hard")
$ws.Range("C2:C43").Validation.Delete()
$ws.Range("C2:C43").Validation.Add(3, 0, 1, "Yes, No")
$ws.Range("D2:D43").Validation.Delete()
$ws.Range("D2:D43").Validation.Add(3, 1, 1, "Yes, No")
$ws.Range("B2:B43").Validation.Delete()
$ws.Range("B2:B43").Validation.Add(3, 1, 1, "Array, Binary, Dynamic Programming, Graph, Interval, Linked List, Matrix, String, Tree, Heap, Class Design")

$ws.Range("J27").Select()

$wb.Save()
